$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for task "1.11" right before the US#2 header row (row 12) ---
$ws.Rows("12").Insert()

# Grow the merged "N. US" / "Descricao" cells for US#1 (rows 2-11 -> 2-12)
# so the new task row still belongs visually to user story #1.
$ws.Range("A2:A12").Merge()
$ws.Range("B2:B12").Merge()

# --- Fill in the new task row (1.11) ---
# D (task number) typed first, then C (description) so the shared-string table
# ends up with "1.11" before the long description, matching authoring order.
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.11"
$ws.Range("C12").Value = "Criar estrutura inicial do backend com repositórios genéricos e unitOfWork e métodos assincronos"
$ws.Range("E12").Value = "Fazendo"
$ws.Rows("12").RowHeight = 30

# New row is "in progress" -> copy the yellow highlight formatting that used to
# live on the old task 1.9 row (row 10, before the insert) onto the new row.
$ws.Range("C10:E10").Copy()
$ws.Range("C12:E12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Task 1.9 (row 10) is now finished -> restyle it like the other "Concluido"
# rows (copy formatting from row 11) and update its status text.
$ws.Range("C11:E11").Copy()
$ws.Range("C10:E10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("E10").Value = "Concluído"

# --- Selection moves to the cell that now holds "Desenvolver tela de login" ---
$ws.Range("C13").Select()
